$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows appended after the last existing row (519), matching the
# vehicle service log entries added for 2022-10-13 (serial 44847), plus a
# single entry on 2022-10-12 that carries over row numbering.

$ws.Cells.Item(520, 1).Value = 44847
$ws.Cells.Item(520, 2).Value = "KA51MD6026"
$ws.Cells.Item(520, 3).Value = "VERNA"
$ws.Cells.Item(520, 4).Value = "PMS"
$ws.Cells.Item(520, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(520, 6).Value = 3442
$ws.Cells.Item(520, 7).Value = "CREDIT"

$ws.Cells.Item(521, 1).Value = 44847
$ws.Cells.Item(521, 2).Value = "TN10AX6476"
$ws.Cells.Item(521, 3).Value = "H CITY"
$ws.Cells.Item(521, 4).Value = "PMS"
$ws.Cells.Item(521, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(521, 6).Value = 6851
$ws.Cells.Item(521, 7).Value = "CARD"

$ws.Cells.Item(522, 1).Value = 44847
$ws.Cells.Item(522, 2).Value = "KA03AD4804 "
$ws.Cells.Item(522, 3).Value = "VERITO"
$ws.Cells.Item(522, 4).Value = "PMS"
$ws.Cells.Item(522, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(522, 6).Value = 5566
$ws.Cells.Item(522, 7).Value = "PAYTM"

$ws.Cells.Item(523, 1).Value = 44847
$ws.Cells.Item(523, 2).Value = "KA03MZ3385"
$ws.Cells.Item(523, 3).Value = "CRYSTA"
$ws.Cells.Item(523, 4).Value = "BRAKE PAD CHANGE"
$ws.Cells.Item(523, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(523, 6).Value = 2200
$ws.Cells.Item(523, 7).Value = "PAYTM"

$ws.Cells.Item(524, 1).Value = 44847
$ws.Cells.Item(524, 2).Value = "KA03MS1238"
$ws.Cells.Item(524, 3).Value = "FIGO"
$ws.Cells.Item(524, 4).Value = "RUNNING REPAIR"
$ws.Cells.Item(524, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(524, 6).Value = 7238

$ws.Cells.Item(525, 1).Value = 44847
$ws.Cells.Item(525, 2).Value = "KA01MG9760"
$ws.Cells.Item(525, 3).Value = "VISTA"
$ws.Cells.Item(525, 4).Value = "GENERAL CHECKUP"
$ws.Cells.Item(525, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(525, 6).Value = 4642
$ws.Cells.Item(525, 7).Value = "P PAY"

# Leave the cursor where data entry ended, one row below/right of the data.
$ws.Range("H521").Select() | Out-Null
